# Apply the "May 9th" data update to the walkingToRunning sample sheet.
# 3 new sensor-reading rows are inserted right after the header (pushing the
# existing readings down), and 7 new rows are appended at the end, extending
# the used range from A1:H21 to A1:H31. The "timestamp"/"label" columns follow
# the existing deterministic pattern (timestamp = (row-2)*100, label = the
# constant activity name) for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$label = "walkingToRunning"

# ax, ay, az, gx, gy, gz readings for every data row (row 2 = timestamp 0) in
# the final, post-edit layout.
$rows = @(
    @(-9.328194618225098, -6.853240489959717, 5.245095252990723, 0.3769224858049733, -0.2748469442259254, 0.004896030350543701),
    @(-7.405017375946045, -14.76152610778809, 0.4900901317596435, -0.8251284312733457, 0.638507986693405, 1.034206660098391),
    @(-8.944772720336914, -6.256664276123047, 6.951004505157471, -0.384797418221631, 2.191421729508453, 1.820474159769618),
    @(-9.14015007019043, -3.555113315582275, 2.792672872543335, -0.8565583720347749, 1.682360221204802, 1.74579396466501),
    @(-0.2594857215881347, -3.807053089141846, 8.088808059692383, -0.239309732830677, 0.5754852927407828, 1.096922339830855),
    @(4.089113235473633, -6.850908756256104, 16.61727523803711, 0.1751386040748463, -0.2145548291081436, 0.8761335207905833),
    @(6.374373435974121, 6.74084997177124, 8.798618316650391, -1.060679824819472, 0.4943530759957171, 1.891958522744575),
    @(-18.19405746459961, -1.922885417938232, 11.5517110824585, -0.450785677701105, -0.1085982450231191, 0.1413266023694169),
    @(-16.09731292724609, -5.246491432189941, 1.230105638504028, 0.3425711914571603, -2.118668488800253, -2.763528636449275),
    @(1.167733669281006, -27.47990798950196, 20.13128471374512, 1.661536096485417, -2.177148738282218, -3.591206591202183),
    @(38.20926666259766, 8.221505165100098, 37.43196487426758, 1.500267133962643, -2.098495937851318, -5.773556591641946),
    @(-9.185737609863279, 6.29133415222168, 1.174301743507385, -2.984340389743152, -1.351582744756633, 1.246524475547487),
    @(-21.03300285339355, -27.21039581298828, -3.327969312667847, -5.330522572629836, 1.363846176055872, 4.617774886335326),
    @(-0.1153392791748046, 1.662713766098022, -6.602585792541504, -4.297088194101645, 7.46291283436744, 1.662266135736299),
    @(-30.58404541015625, -11.14975070953369, -8.237148284912109, -0.1566145409663497, 4.606385328884262, -0.8842834989056194),
    @(7.598236083984375, -49.00935745239258, 8.44972038269043, 5.259676156606183, -3.144118193455617, -3.28624562092742),
    @(-6.318488121032715, 9.492542266845703, 1.218106150627136, 4.500417068015016, -0.31324529439601, -3.402522266171384),
    @(-51.58052062988281, -63.86045837402344, 42.37603759765625, -1.445235498086886, 2.503525588710229, 0.09940216218542508),
    @(-8.470245361328125, 5.231861114501953, -25.23441886901855, -8.082223704808689, 0.7967113996306803, 8.42470887134173),
    @(-21.50184631347656, -2.417125701904297, -8.449828147888184, -1.0955425516487, -6.407353537588252, 3.211009324377832),
    @(-21.86991500854492, -24.83074569702148, 1.575617074966431, 3.420750588829331, -8.409929646154669, 0.8056608922617063),
    @(3.400870323181152, 13.92329216003418, 5.784989833831787, 5.383590975182555, -4.346148073412959, -0.09347351565632778),
    @(0.5669898986816406, 13.84309577941894, 9.06630802154541, 1.778205815361046, -1.402306104331036, -0.230236690637843),
    @(4.273637771606445, -40.02112579345703, 19.01494979858398, -3.16472590750482, 0.5638766811924731, 1.501450692722371),
    @(-33.29051208496094, 13.78821468353272, -5.866414546966553, -2.917193829232435, 0.08571250834320598, -0.5273152324309833),
    @(-38.59527587890625, -11.46046161651611, -14.89169502258301, -0.5172856062260307, 4.862342070283832, -2.60658347866941),
    @(69.93356323242188, -70.74887084960938, 44.17119979858398, 1.666729755276673, 0.2906178595197595, -1.515296798606113),
    @(-14.73312759399414, 12.10027313232422, 5.143325805664063, 2.873153332539526, -9.141991994266446, -1.962461360677373),
    @(16.71687316894531, -16.55658340454102, 22.30324172973633, 1.324421677787197, -4.62868054152566, -0.05621947888198298),
    @(-58.74892807006836, -18.46404838562012, -22.14792823791504, -3.513774491813999, 3.990525447645534, 6.551176387670226)
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $i * 100
    $ws.Cells.Item($r, 2).Value = $label
    $vals = $rows[$i]
    $ws.Cells.Item($r, 3).Value = $vals[0]
    $ws.Cells.Item($r, 4).Value = $vals[1]
    $ws.Cells.Item($r, 5).Value = $vals[2]
    $ws.Cells.Item($r, 6).Value = $vals[3]
    $ws.Cells.Item($r, 7).Value = $vals[4]
    $ws.Cells.Item($r, 8).Value = $vals[5]
}

